# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (with the same per-quarter fund-holdings
# layout as the existing 2021-Qx sheets) right before the "总计" (totals)
# summary sheet, populates it with the 2022-Q1 holdings, and inserts a new
# leading row into the "总计" sheet summarising that quarter.

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2021-Q3")
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet immediately before "总计"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($total)
$ws.Name = "2022-Q1"

# Reuse the existing header / index-column formatting from sibling sheets
# so no spurious new styles get created.
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("A2:A11").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Holdings rows (A = index, B..G = text-typed, H = numeric rank)
# ---------------------------------------------------------------------
$rows = @(
    @(0, "005777", "广发科技动力股票",                       "17.47", "92.65", "7.08", "1.2369", 3),
    @(1, "213008", "宝盈资源优选混合",                         "11.31", "81.52", "4.04", "0.4569", 6),
    @(2, "007731", "民生加银持续成长混合A",                   "2.62",  "93.83", "7.29", "0.1910", 1),
    @(3, "004314", "前海开源沪港深新硬件主题灵活配置混合A",   "1.67",  "90.05", "4.32", "0.0721", 9),
    @(4, "002707", "摩根士丹利华鑫科技领先灵活配置混合",     "2.27",  "93.05", "3.09", "0.0701", 10),
    @(5, "004315", "前海开源沪港深新硬件主题灵活配置混合C",   "1.00",  "90.05", "4.32", "0.0432", 9),
    @(6, "004044", "金鹰转型动力灵活配置混合",                 "0.72",  "93.34", "3.35", "0.0241", 10),
    @(7, "007732", "民生加银持续成长混合C",                   "0.14",  "93.83", "7.29", "0.0102", 1),
    @(8, "001914", "中信建投聚利混合A",                       "0.13",  "39.07", "2.07", "0.0027", 7),
    @(9, "000041", "华夏全球精选股票(QDII)",                  "0.02",  "39.07", "2.07", "0.0004", 7)
)

# Force columns B..G to be stored as text (so codes keep leading zeros and
# numbers like "17.47" stay text, matching the source data), then drop the
# number-format override again once the values are in so no extra style
# sticks around on the cells.
$dataRange = $ws.Range("B2:G11")
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

$dataRange.ClearFormats()

# ---------------------------------------------------------------------
# 4. Insert the 2022-Q1 summary row at the top of "总计", pushing the
#    existing quarters down one row.
# ---------------------------------------------------------------------
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$b4 = $total.Range("B4").Value()
$c4 = $total.Range("C4").Value()
$d4 = $total.Range("D4").Value()

$b3 = $total.Range("B3").Value()
$c3 = $total.Range("C3").Value()
$d3 = $total.Range("D3").Value()

$b2 = $total.Range("B2").Value()
$c2 = $total.Range("C2").Value()
$d2 = $total.Range("D2").Value()

$total.Range("B5").Value = $b4
$total.Range("C5").Value = $c4
$total.Range("D5").Value = $d4

$total.Range("B4").Value = $b3
$total.Range("C4").Value = $c3
$total.Range("D4").Value = $d3

$total.Range("B3").Value = $b2
$total.Range("C3").Value = $c2
$total.Range("D3").Value = $d2

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 2.11

Write-Output "2022-Q1 sheet added and 总计 updated"
